$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competição") values were incorrectly truncated from 270 to 70.
# Restore the dropped leading digit for every data row (rows 2-66).
$ws.Range("B2:B66").Value = 270
